$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 46
$ws.Range("F2").Value = 30
$ws.Range("D3").Value = 25
$ws.Range("F3").Value = 47
$ws.Range("D4").Value = 44
$ws.Range("F4").Value = 26
